# Add a new sale line ("كريم 555") to the DaySale report as row 46,
# pushing the existing totals/footer rows down by one, updating the
# grand total and refreshing the generated-on timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a fresh row above the "Total" row (old row 46) ------------
$ws.Rows("46:46").Insert()

# Clone the formatting of the last existing data row (45) into the new
# row 46 so fonts/fills/number-formats/merges match the other product
# rows exactly.
$ws.Range("A45:Q45").Copy()
$ws.Range("A46:Q46").PasteSpecial(-4104)   # xlPasteAll

# The thin row-separator (bottom) border sometimes gets dropped by the
# insert/paste - reapply it explicitly so row 46 matches its neighbours.
$newRow = $ws.Range("A46:Q46")
$newRow.Borders.Item(9).LineStyle = 1
$newRow.Borders.Item(9).Color = 13882323

# Match row height used by the other product rows.
$ws.Rows("46:46").RowHeight = 25.5

# --- Populate the new row's values -------------------------------------
$ws.Range("A46").Value2 = 40
$ws.Range("C46").Value2 = "كريم 555"
$ws.Range("H46").Value2 = "3:0"

# Columns L and P store their numbers as text in this workbook, so force
# a text number-format while assigning, then restore the original
# number-format (keeps the same style as the rest of the sheet).
$lFmt = $ws.Range("L46").NumberFormat
$ws.Range("L46").NumberFormat = "@"
$ws.Range("L46").Value2 = "0"
$ws.Range("L46").NumberFormat = $lFmt

$ws.Range("N46").Value2 = "20.00"

$pFmt = $ws.Range("P46").NumberFormat
$ws.Range("P46").NumberFormat = "@"
$ws.Range("P46").Value2 = "20.0000"
$ws.Range("P46").NumberFormat = $pFmt

$ws.Range("Q46").Value2 = "1:0"

# --- Update the totals row (now row 47) --------------------------------
$ws.Range("P47").Value2 = 1518.8

# --- Refresh the "generated on" timestamp in the footer (now row 48) --
$ws.Range("A48").Value2 = "Wednesday, 16 July, 2025 3:24 PM"
